# Updated ERD Image file.
# Rename column/field header labels from Title_Case to lower_snake_case
# to match the updated ERD image conventions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Categories table header column
$ws.Range("B2").Value = "name"

# Restaurants table headers
$ws.Range("B8").Value = "name"
$ws.Range("C8").Value = "street"
$ws.Range("D8").Value = "city"
$ws.Range("E8").Value = "state"
$ws.Range("F8").Value = "zip_code"
$ws.Range("G8").Value = "category_id"

# Users table headers
$ws.Range("B14").Value = "first_name"
$ws.Range("C14").Value = "last_name"

# Reviews table headers
$ws.Range("B19").Value = "comment"
$ws.Range("C19").Value = "rating"
$ws.Range("D19").Value = "restaurant_id"
$ws.Range("E19").Value = "user_id"

# Update the active selection to match the saved view state
$ws.Range("D18").Select()
